$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the customer name in B2 from "tanmay" to "Tanmay Bore"
$ws.Range("B2").Value = "Tanmay Bore"
